$d = $word.ActiveDocument

# 1. Trim the leading fragment off the very first run of the document:
#    "даленный репозиторий, ... команды могли увидеть их." -> "ды могли увидеть их."
$trimRange = $d.Range(0, $d.Content.End)
$trimRange.Find.Execute(
    "даленный репозиторий, чтобы другие участники коман", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 1)

# 2. Duplicate the "Без смысла" + explanatory sentence block and insert
#    the copy right before the (now shifted) first occurrence of that
#    block, so it reads: "...их.Без смысла <sentence>Без смысла <sentence>..."
$findRange = $d.Range(0, $d.Content.End)
$findRange.Find.Execute("Без смысла", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertionPoint = $d.Range($findRange.Start, $findRange.Start)
$insertionPoint.InsertBefore("Без смысла используется для отправки ваших локальных изменений (коммитов) в удаленный репозиторий. Это означает, что все ваши локальные изменения будут скопированы в удаленный репозиторий, чтобы другие участники команды могли увидеть их.")
